$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 115: 2024-07-29 data point for PAL.MI
$ws.Range("A115").Value = 45502.2916666667
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("B115").Value = 0
$ws.Range("C115").Value = 5.94000005722046
$ws.Range("D115").Value = 5.94000005722046
$ws.Range("E115").Value = 5.94000005722046
$ws.Range("F115").Value = 5.94000005722046
$ws.Range("G115").NumberFormat = "@"
$ws.Range("G115").Value = "5.94000005722046"
$ws.Range("G115").ClearFormats()
$ws.Range("H115").Value = "PAL.MI"

# New row 116: 2024-07-30 data point for PAL.MI
$ws.Range("A116").Value = 45503.6436805556
$ws.Range("A114").Copy()
$ws.Range("A116").PasteSpecial(-4122)
$ws.Range("B116").Value = 1200
$ws.Range("C116").Value = 6.03999996185303
$ws.Range("D116").Value = 6
$ws.Range("E116").Value = 6.01999998092651
$ws.Range("F116").Value = 6
$ws.Range("G116").NumberFormat = "@"
$ws.Range("G116").Value = "6"
$ws.Range("G116").ClearFormats()
$ws.Range("H116").Value = "PAL.MI"

$excel.CutCopyMode = $false
